$d = $word.ActiveDocument

# --- Change 1: expand algorithm list with approach descriptions ---
$d.Content.Find.Execute(
    "Naïve Bayes, Support Vector Machines and Random Forest",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Naïve Bayes (Probabilistic approach), Support Vector Machines (Hyperplane approach) and Random Forest",
    2) | Out-Null

# --- Change 3: add vocab size after first min_df 5 mention ---
$d.Content.Find.Execute(
    "Unigrams & Bigrams along with min_df 5 gave both the best features along with",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Unigrams & Bigrams along with min_df 5(Vocab Size 56619) gave both the best features along with",
    2) | Out-Null

# --- Change 4: add vocab size after second min_df of 5 mention ---
$d.Content.Find.Execute(
    "min_df of 5 capturing words with higher weights",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "min_df of 5(Vocab Size 56619) capturing words with higher weights",
    2) | Out-Null

# --- Change 5: add vocab size after Min_df=5 mention ---
$d.Content.Find.Execute(
    "Min_df=5 gave the words and performance better than higher min_df",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Min_df=5(Vocab Size 72390) gave the words and performance better than higher min_df",
    2) | Out-Null

# --- Change 6: simplify Random Forest sentence ---
$d.Content.Find.Execute(
    "Random Forest Confusion matrix and classification report showed a very low",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Random Forest showed a very low",
    2) | Out-Null

# --- Change 7: drop ", eat healthier" from the common-words list ---
$d.Content.Find.Execute(
    "lose, weight, diet, eat healthier which are used",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "lose, weight, diet which are used",
    2) | Out-Null

# --- Move the _GoBack bookmark to its new location ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$rng = $d.Content
$rng.Find.Execute("lose, weight, diet ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$target = $d.Range($rng.End, $rng.End)
$d.Bookmarks.Add("_GoBack", $target) | Out-Null
